$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1 with the same style as the other header cells (copy formatting from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for I2:J82
$ijData = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(9, 9)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(10, 10)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(10, 10)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(9, 9)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(7, 7)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(10, 10)
    35 = @(8, 8)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(9, 9)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(8, 9)
    44 = @(8, 8)
    45 = @(9, 9)
    46 = @(9, 9)
    47 = @(9, 9)
    48 = @(9, 9)
    49 = @(9, 9)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(9, 9)
    53 = @(8, 8)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(7, 7)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(9, 9)
    67 = @(9, 9)
    68 = @(9, 9)
    69 = @(9, 9)
    70 = @(9, 9)
    71 = @(9, 9)
    72 = @(9, 9)
    73 = @(9, 9)
    74 = @(8, 8)
    75 = @(9, 9)
    76 = @(8, 8)
    77 = @(6, 6)
    78 = @(6, 6)
    79 = @(6, 6)
    80 = @(6, 6)
    81 = @(6, 6)
    82 = @(4, 4)
}

foreach ($r in $ijData.Keys) {
    $vals = $ijData[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
